# Adyen test fixture migration:
#   - "Gross Currency" (col K) / "Net Currency" (col O) values of "EUR"
#     are updated to "USD" for every data row.
#   - Selection / active cell is moved to P38 to match the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 5
$lastDataRow  = 44
$currencyCols = @("K", "O")

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    foreach ($col in $currencyCols) {
        $cell = $ws.Range($col + $row)
        if ($cell.Value2 -eq "EUR") {
            $cell.Value = "USD"
        }
    }
}

$null = $ws.Range("P38").Select()
